$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: stage original row contents into scratch rows (100+) to avoid overwrite-before-read issues
$ws.Range("A2:AY2").Copy($ws.Range("A100"))
$ws.Range("A3:AY3").Copy($ws.Range("A101"))
$ws.Range("A4:AY4").Copy($ws.Range("A102"))
$ws.Range("A5:AY5").Copy($ws.Range("A103"))
$ws.Range("A6:AY6").Copy($ws.Range("A104"))
$ws.Range("A7:AY7").Copy($ws.Range("A105"))
$ws.Range("A8:AY8").Copy($ws.Range("A106"))
$ws.Range("A10:AY10").Copy($ws.Range("A107"))
$ws.Range("A11:AY11").Copy($ws.Range("A108"))
$ws.Range("A13:AY13").Copy($ws.Range("A109"))
$ws.Range("A14:AY14").Copy($ws.Range("A110"))
$ws.Range("A15:AY15").Copy($ws.Range("A111"))
$ws.Range("A16:AY16").Copy($ws.Range("A112"))
$ws.Range("A17:AY17").Copy($ws.Range("A113"))

# Step 2: write staged content into destination rows per the permutation
$ws.Range("A102:AY102").Copy($ws.Range("A2"))
$ws.Range("A105:AY105").Copy($ws.Range("A3"))
$ws.Range("A100:AY100").Copy($ws.Range("A4"))
$ws.Range("A104:AY104").Copy($ws.Range("A5"))
$ws.Range("A111:AY111").Copy($ws.Range("A6"))
$ws.Range("A107:AY107").Copy($ws.Range("A7"))
$ws.Range("A110:AY110").Copy($ws.Range("A8"))
$ws.Range("A101:AY101").Copy($ws.Range("A10"))
$ws.Range("A113:AY113").Copy($ws.Range("A11"))
$ws.Range("A112:AY112").Copy($ws.Range("A13"))
$ws.Range("A108:AY108").Copy($ws.Range("A14"))
$ws.Range("A103:AY103").Copy($ws.Range("A15"))
$ws.Range("A106:AY106").Copy($ws.Range("A16"))
$ws.Range("A109:AY109").Copy($ws.Range("A17"))

# Step 3: clear scratch rows
$ws.Range("A100:AY100").ClearContents()
$ws.Range("A101:AY101").ClearContents()
$ws.Range("A102:AY102").ClearContents()
$ws.Range("A103:AY103").ClearContents()
$ws.Range("A104:AY104").ClearContents()
$ws.Range("A105:AY105").ClearContents()
$ws.Range("A106:AY106").ClearContents()
$ws.Range("A107:AY107").ClearContents()
$ws.Range("A108:AY108").ClearContents()
$ws.Range("A109:AY109").ClearContents()
$ws.Range("A110:AY110").ClearContents()
$ws.Range("A111:AY111").ClearContents()
$ws.Range("A112:AY112").ClearContents()
$ws.Range("A113:AY113").ClearContents()
